$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    # Force the cell to remain plain text even when the new
    # value looks numeric (e.g. "0.657"), then restore the
    # cells original style so no formatting changes leak in.
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "43.708.04"
$ws.Range("E2").Value = "  +0.80%  "
Set-TextValue $ws "D3" "2.349.73"
$ws.Range("E3").Value = "  +4.41%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws "D5" "234.97"
$ws.Range("E5").Value = "  +1.87%  "
Set-TextValue $ws "D6" "0.657"
$ws.Range("E6").Value = "  +2.80%  "
Set-TextValue $ws "D7" "73.42"
$ws.Range("E7").Value = "  +14.04%  "
$ws.Range("E8").Value = "  -0.01%  "
Set-TextValue $ws "D9" "0.526"
$ws.Range("E9").Value = "  +20.19%  "
Set-TextValue $ws "D10" "0.0981"
$ws.Range("E10").Value = "  +3.40%  "
Set-TextValue $ws "D11" "27.48"
$ws.Range("E11").Value = "  +4.48%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws "D12" "0.106"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws "D13" "17.06"
$ws.Range("E13").Value = "  +14.79%  "
Set-TextValue $ws "D14" "2.698.75"
$ws.Range("E14").Value = "  +4.46%  "
Set-TextValue $ws "D15" "6.68"
$ws.Range("E15").Value = "  +10.84%  "
Set-TextValue $ws "D16" "0.881"
$ws.Range("E16").Value = "  +7.64%  "
Set-TextValue $ws "D17" "2.345.37"
$ws.Range("E17").Value = "  +4.15%  "
Set-TextValue $ws "D18" "43.548.31"
$ws.Range("E18").Value = "  +0.73%  "
Set-TextValue $ws "D19" "0.0000100"
$ws.Range("E19").Value = "  +4.13%  "
Set-TextValue $ws "D20" "76.10"
$ws.Range("E20").Value = "  +4.44%  "
Set-TextValue $ws "D21" "6.36"
$ws.Range("E21").Value = "  +5.12%  "
Set-TextValue $ws "D22" "250.39"
$ws.Range("E22").Value = "  +1.46%  "
Set-TextValue $ws "D23" "3.82"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("E24").Value = "  +0.01%  "
Set-TextValue $ws "D25" "2.48"
$ws.Range("E25").Value = "  +3.21%  "
Set-TextValue $ws "D26" "10.28"
$ws.Range("E26").Value = "  +6.14%  "
Set-TextValue $ws "D27" "2.24"
$ws.Range("E27").Value = "  -2.13%  "
Set-TextValue $ws "D28" "22.39"
$ws.Range("E28").Value = "  +3.90%  "
Set-TextValue $ws "D29" "172.30"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("E30").Value = "  +7.88%  "
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("E32").Value = "  +4.33%  "
$ws.Range("E33").Value = "  +3.85%  "
Set-TextValue $ws "D34" "0.0699"
$ws.Range("E34").Value = "  +3.29%  "
Set-TextValue $ws "D35" "5.10"
$ws.Range("E35").Value = "  +3.98%  "
Set-TextValue $ws "D36" "3.77"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("E37").Value = "  +7.45%  "
Set-TextValue $ws "D38" "6.40"
$ws.Range("E38").Value = "  +0.42%  "
Set-TextValue $ws "D39" "0.0265"
$ws.Range("E39").Value = "  +6.43%  "
Set-TextValue $ws "D40" "19.45"
$ws.Range("E40").Value = "  +13.83%  "
$ws.Range("E41").Value = "  +0.00%  "
Set-TextValue $ws "D42" "8.89"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("E43").Value = "  +8.60%  "
Set-TextValue $ws "D44" "1.21"
$ws.Range("E44").Value = "  +3.24%  "
Set-TextValue $ws "D45" "98.66"
$ws.Range("E45").Value = "  +2.35%  "
Set-TextValue $ws "D46" "0.0962"
$ws.Range("E46").Value = "  +2.82%  "
Set-TextValue $ws "D47" "4.44"
$ws.Range("E47").Value = "  -0.41%  "
Set-TextValue $ws "D48" "0.179"
$ws.Range("E48").Value = "  +12.88%  "
Set-TextValue $ws "D49" "1.438.82"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  +2.33%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws "D51" "2.571.83"
$ws.Range("E51").Value = "  +4.43%  "
